# Increase substitution potential to 20% in 2022 and 50% in 2030
# Scen_B_IND_Mitigations.xlsx -- sheet "IND_Shares"
#
# Across the mitigation table, the "2022" substitution-potential column (R)
# goes from 5% to 20%, and the "2030" substitution-potential column (S) goes
# from 30% to 50%, for every mitigation row that carried the old 5%/30%
# pair. Two rows (84 and 90) encode the same 5%/30% -> 20%/50% change via
# percentage formulas instead of plain literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IND_Shares")
$ws.Activate()

# Rows whose R (2022) / S (2030) substitution-potential cells are plain
# literals 0.05 / 0.3 that become 0.2 / 0.5.
$rows = @(12, 13, 14, 24, 26, 27, 34, 35, 36, 44, 45, 46, 53, 54, 55, 64, 65, 66, 73, 74, 75, 87, 88, 89, 97, 98, 99, 106, 107, 108, 115, 116, 117, 124, 125, 126, 134, 135, 136, 143, 144, 145)

foreach ($r in $rows) {
    $ws.Range("R$r").Value = 0.2
    $ws.Range("S$r").Value = 0.5
}

# Row 84: R84 was "=5%*80%" (i.e. 5% potential) -- now a plain 20%*80% value.
# S84 stays a formula, but the percentage bumps from 30% to 50%.
$ws.Range("R84").Value = 0.16
$ws.Range("S84").Formula = "=50%*80%"

# Row 90: only S90's formula potential bumps from 30% to 50% (R90 keeps 5%).
$ws.Range("S90").Formula = "=50%*20%"

# Restore the view to the default top-left / A1 selection (the saved file
# had scrolled to A34 with T72 selected).
$ws.Range("A1").Select()
